$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append new ticker row at the end (row 65), matching the existing data pattern
$ws.Range("A65").Value = "GRT-USD"
